$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-65 down to 13-66
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with new data
$ws.Range("A12").Value2 = 5
$ws.Range("B12").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C12").Value2 = 'Maule'
$ws.Range("D12").Value2 = 44831
$ws.Range("E12").Value2 = 7
$ws.Range("F12").Value2 = 300000000
$ws.Range("G12").Value2 = 'Espárragos'
$ws.Range("H12").Value2 = 'Sin especificar'
$ws.Range("I12").Value2 = 'Primera'
$ws.Range("J12").Value2 = 1000
$ws.Range("K12").Value2 = 2000
$ws.Range("L12").Value2 = 2000
$ws.Range("M12").Value2 = 2000
$ws.Range("N12").Value2 = '$/kilo'
$ws.Range("O12").Value2 = 'Provincia de Linares'
$ws.Range("P12").Value2 = 2000
$ws.Range("Q12").Value2 = 1
$ws.Range("R12").Value2 = 'Hortaliza'

# Append new row 67 with new data (after the shift, existing data occupies rows 1-66)
$ws.Range("A67").Value2 = 5
$ws.Range("B67").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C67").Value2 = 'Maule'
$ws.Range("D67").Value2 = 44832
$ws.Range("E67").Value2 = 7
$ws.Range("F67").Value2 = 300000000
$ws.Range("G67").Value2 = 'Espárragos'
$ws.Range("H67").Value2 = 'Sin especificar'
$ws.Range("I67").Value2 = 'Primera'
$ws.Range("J67").Value2 = 1500
$ws.Range("K67").Value2 = 2000
$ws.Range("L67").Value2 = 2000
$ws.Range("M67").Value2 = 2000
$ws.Range("N67").Value2 = '$/kilo'
$ws.Range("O67").Value2 = 'Provincia de Linares'
$ws.Range("P67").Value2 = 2000
$ws.Range("Q67").Value2 = 1
$ws.Range("R67").Value2 = 'Hortaliza'

# Row 67 is a brand new row outside the previous used range, so it does not
# automatically inherit the date number format used by column D elsewhere;
# set it explicitly to match the rest of the column.
$ws.Range("D67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
